$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''26.246.03'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '''  +0.30%  '
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = '''1.605.01'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '''  +0.07%  '
$ws.Range("E3").ClearFormats()
$ws.Range("E4").Value = '''  -0.10%  '
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = '''212.53'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '''  -0.13%  '
$ws.Range("E5").ClearFormats()
$ws.Range("E6").Value = '''  -0.12%  '
$ws.Range("E6").ClearFormats()
$ws.Range("D8").Value = '''0.250'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '''  +0.59%  '
$ws.Range("E8").ClearFormats()
$ws.Range("E9").Value = '''  -0.33%  '
$ws.Range("E9").ClearFormats()
$ws.Range("D10").Value = '''18.36'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '''  +1.98%  '
$ws.Range("E10").ClearFormats()
$ws.Range("D11").Value = '''0.0813'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '''  -0.72%  '
$ws.Range("E11").ClearFormats()
$ws.Range("E12").Value = '''  +0.04%  '
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = '''1.619.55'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '''  +1.05%  '
$ws.Range("E13").ClearFormats()
$ws.Range("E14").Value = '''  +0.46%  '
$ws.Range("E14").ClearFormats()
$ws.Range("D15").Value = '''0.514'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '''  +0.58%  '
$ws.Range("E15").ClearFormats()
$ws.Range("D16").Value = '''26.211.16'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '''  +0.20%  '
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = '''61.95'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '''  +2.47%  '
$ws.Range("E17").ClearFormats()
$ws.Range("D18").Value = '''0.0₃0728'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '''  +0.73%  '
$ws.Range("E18").ClearFormats()
$ws.Range("E19").Value = '''  -0.07%  '
$ws.Range("E19").ClearFormats()
$ws.Range("D20").Value = '''200.07'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '''  -2.34%  '
$ws.Range("E20").ClearFormats()
$ws.Range("E21").Value = '''  +0.53%  '
$ws.Range("E21").ClearFormats()
$ws.Range("E22").Value = '''  +0.07%  '
$ws.Range("E22").ClearFormats()
$ws.Range("D23").Value = '''6.01'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '''  +0.32%  '
$ws.Range("E23").ClearFormats()
$ws.Range("E24").Value = '''  +2.84%  '
$ws.Range("E24").ClearFormats()
$ws.Range("D25").Value = '''143.96'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '''  +1.72%  '
$ws.Range("E25").ClearFormats()
$ws.Range("E26").Value = '''  -0.06%  '
$ws.Range("E26").ClearFormats()
$ws.Range("E27").Value = '''  -2.27%  '
$ws.Range("E27").ClearFormats()
$ws.Range("D28").Value = '''15.18'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '''  -0.05%  '
$ws.Range("E28").ClearFormats()
$ws.Range("E29").Value = '''  +2.00%  '
$ws.Range("E29").ClearFormats()
$ws.Range("D30").Value = '''0.0492'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '''  +4.55%  '
$ws.Range("E30").ClearFormats()
$ws.Range("E31").Value = '''  +0.46%  '
$ws.Range("E31").ClearFormats()
$ws.Range("E32").Value = '''  +2.75%  '
$ws.Range("E32").ClearFormats()
$ws.Range("D33").Value = '''2.95'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '''  -1.33%  '
$ws.Range("E33").ClearFormats()
$ws.Range("E34").Value = '''  +0.95%  '
$ws.Range("E34").ClearFormats()
$ws.Range("D35").Value = '''2.37'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '''  +1.04%  '
$ws.Range("E35").ClearFormats()
$ws.Range("D36").Value = '''1.165.30'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '''  +4.30%  '
$ws.Range("E36").ClearFormats()
$ws.Range("D37").Value = '''0.0169'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '''  +3.17%  '
$ws.Range("E37").ClearFormats()
$ws.Range("D39").Value = '''2.32'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '''  +1.07%  '
$ws.Range("E39").ClearFormats()
$ws.Range("E40").Value = '''  +0.10%  '
$ws.Range("E40").ClearFormats()
$ws.Range("D41").Value = '''0.496'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '''  +0.92%  '
$ws.Range("E41").ClearFormats()
$ws.Range("E42").Value = '''  +4.13%  '
$ws.Range("E42").ClearFormats()
$ws.Range("E43").Value = '''  +0.40%  '
$ws.Range("E43").ClearFormats()
$ws.Range("D44").Value = '''1.739.45'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '''  +0.04%  '
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = '''91.92'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '''  -1.02%  '
$ws.Range("E45").ClearFormats()
$ws.Range("D46").Value = '''0.0₆0106'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '''  +15.16%  '
$ws.Range("E46").ClearFormats()
$ws.Range("D47").Value = '''1.54'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '''  +1.59%  '
$ws.Range("E47").ClearFormats()
$ws.Range("D48").Value = '''53.99'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '''  +0.95%  '
$ws.Range("E48").ClearFormats()
$ws.Range("E49").Value = '''  +0.12%  '
$ws.Range("E49").ClearFormats()
